$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct small floating point revisions in existing cells
$ws.Range("D9").Value = 20.079509735107422
$ws.Range("D19").Value = 21.123500823974609
$ws.Range("D20").Value = 124.09785461425781
$ws.Range("E20").Value = 119.66958618164063

# Add new row 23 for lccNA_pcnt
$ws.Range("A23").Value = "lccNA_pcnt"
$ws.Range("B23").Value = 36.570178985595703
$ws.Range("C23").Value = 37.286026000976563
$ws.Range("D23").Value = 38.187744140625
$ws.Range("E23").Value = 39.310688018798828
$ws.Range("F23").Value = 40.399639129638672
$ws.Range("G23").Value = 41.062961578369141
$ws.Range("H23").Value = 41.443065643310547

# Apply same style as other numeric rows (style index 1 -> numFmtId 1, border)
$ws.Range("B23:H23").NumberFormat = $ws.Range("B22:H22").NumberFormat
$ws.Range("B23:H23").Borders.LineStyle = $ws.Range("B22:H22").Borders.LineStyle
